$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 296 (shifts old rows 296:367 down to 297:368,
# carrying their formatting/values along for the ride).
$ws.Rows.Item(296).Insert()

# Populate the newly inserted row 296 with the new record.
$ws.Cells.Item(296, 1).Value = 11
$ws.Cells.Item(296, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(296, 3).Value = "Bíobío"
$ws.Cells.Item(296, 4).Value = 45015
$ws.Cells.Item(296, 5).Value = 8
$ws.Cells.Item(296, 6).Value = 100112045
$ws.Cells.Item(296, 7).Value = "Zapallo"
$ws.Cells.Item(296, 8).Value = "Camote"
$ws.Cells.Item(296, 9).Value = "1a (cosecha)"
$ws.Cells.Item(296, 10).Value = 700
$ws.Cells.Item(296, 11).Value = 300
$ws.Cells.Item(296, 12).Value = 350
$ws.Cells.Item(296, 13).Value = 329
$ws.Cells.Item(296, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(296, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(296, 16).Value = 329
$ws.Cells.Item(296, 17).Value = 1
$ws.Cells.Item(296, 18).Value = "Hortaliza"
